$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column G (old "synonym" column shifts from G to H)
$ws.Columns("G").Insert()

# New column G holds the "Wolverton 2007" measurement-code mapping.
# Write the data value before the header so the shared-strings table
# picks up "AST4" ahead of "Wolverton 2007" (matches original authoring order).
$ws.Range("G8").Value = "AST4"
$ws.Range("G1").Value = "Wolverton 2007"

# Leave selection on the newly entered cell, matching the author's workflow
$ws.Range("G8").Select()
